# Insert a new "bioconductor" column into the co-occurrence matrix on Sheet1.
# The sheet already has a "bioconductor" row (row 5); we add the matching
# column between "bioconda" (D) and "biolinux" (old E, now shifts to F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns E:M -> F:N, leaving a blank column E for the new data.
$ws.Columns("E:E").Insert()

# Header for the new column.
$ws.Range("E1").Value2 = "bioconductor"

# New column values (co-occurrence counts vs. "bioconductor"), row by row.
$ws.Range("E2").Value2 = 0
$ws.Range("E3").Value2 = 23
$ws.Range("E4").Value2 = 14
$ws.Range("E5").Value2 = 1823
$ws.Range("E6").Value2 = 2
$ws.Range("E7").Value2 = 12
$ws.Range("E8").Value2 = 0
$ws.Range("E9").Value2 = 17
$ws.Range("E10").Value2 = 14
$ws.Range("E11").Value2 = 196
$ws.Range("E12").Value2 = 6
$ws.Range("E13").Value2 = 0
$ws.Range("E14").Value2 = 3
